$d = $word.ActiveDocument

$replacements = @(
    @("293÷7=", "110÷4="),
    @("331÷7=", "299÷4="),
    @("532÷4=", "632÷2="),
    @("982÷9=", "899÷8="),
    @("980÷7=", "440÷2="),
    @("657÷2=", "709÷8="),
    @("649÷9=", "634÷6="),
    @("226÷7=", "406÷8="),
    @("761÷6=", "806÷4="),
    @("889÷5=", "733÷4="),
    @("409÷5=", "928÷7="),
    @("217÷6=", "348÷5="),
    @("258÷9=", "206÷9="),
    @("672÷9=", "649÷2="),
    @("244÷6=", "689÷8="),
    @("554÷5=", "274÷9="),
    @("321÷3=", "540÷9="),
    @("816÷5=", "513÷9="),
    @("443÷3=", "672÷6="),
    @("260÷4=", "755÷6="),
    @("203÷4=", "481÷4="),
    @("718÷7=", "765÷9="),
    @("443÷6=", "817÷3="),
    @("146÷5=", "601÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
